# Update the row-2 data with the newly added iAuthor TC's values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "jTwOQ444"
$ws.Range("B2").Value = 23102030
$ws.Range("C2").Value = "ngccdjv54"
$ws.Range("D2").Value = 'A&$s4z5T'
$ws.Range("F2").Value = "hFFCTWAM"
$ws.Range("G2").Value = "pTFl"
